$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the IMPEDANCE (column C) values for several land-cover rows ---
$ws.Range("C3").Value2  = 85
$ws.Range("C4").Value2  = 5
$ws.Range("C5").Value2  = 10
$ws.Range("C6").Value2  = 15
$ws.Range("C7").Value2  = 20
$ws.Range("C8").Value2  = 30
$ws.Range("C9").Value2  = 40
$ws.Range("C10").Value2 = 45
$ws.Range("C12").Value2 = 35
$ws.Range("C13").Value2 = 45
$ws.Range("C14").Value2 = 45
$ws.Range("C15").Value2 = 20
$ws.Range("C16").Value2 = 45
$ws.Range("C18").Value2 = 25
$ws.Range("C19").Value2 = 25

# --- Row 1 header: taller row, drop the stray column D header cell ---
$ws.Rows.Item(1).RowHeight = 45

# --- Row 27 description wraps onto two lines, so give it more height ---
$ws.Rows.Item(27).RowHeight = 30
$ws.Range("B27").WrapText = $true

# --- Resize the remaining two data columns ---
$ws.Columns.Item(2).ColumnWidth = 28.16666667
$ws.Columns.Item(3).ColumnWidth = 10.88

# --- Column D is no longer part of the table; remove it entirely ---
$ws.Columns.Item(4).Delete()

# --- The database range no longer includes column D ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Database") {
        $n.RefersTo = "=LandCover_Class!`$A`$1:`$C`$30"
    }
}

# --- Selection moves to the (now blank) column D downward ---
$ws.Range("D1:K1048576").Select()
